$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing test data but keep cell formatting intact (so C5 keeps its
# original font style even though it becomes empty again).
$ws.Range("A1:C6").ClearContents()

# --- Header row ---
$ws.Range("A1").Value = "tradeName"
$ws.Range("B1").Value = "Comment"

# A3 is a genuine number (3811989) that simply has the Text format applied on
# top of it afterwards, so it must be entered as a real number *before* any
# Text format is applied to it, otherwise Excel would store it as text like
# the rest of the column.
$ws.Range("A3").Value = 3811989

# --- Apply the "Text" number format (and left alignment for column A rows
#     2-6) to the new negative test-case cells *before* typing values into
#     the remaining ones, so strings that look numeric (e.g. "000000") are
#     kept as text instead of being auto-coerced to numbers.
$ws.Range("A2:A6").NumberFormat = "@"
$ws.Range("B2:B9").NumberFormat = "@"
$ws.Range("A7:A9").NumberFormat = "@"
$ws.Range("A2:A6").HorizontalAlignment = -4131

# --- New negative test-case rows ---
$ws.Range("A2").Value = "Glazing Solubtions"
$ws.Range("B2").Value = "typo"

$ws.Range("B3").Value = "above company extra digi"

$ws.Range("A4").Value = "z1"
$ws.Range("B4").Value = "combi of char and dig"

$ws.Range("A5").Value = "000000"
$ws.Range("B5").Value = "6 zeros"

$ws.Range("A6").Value = "oooooo"
$ws.Range("B6").Value = "6 lower case o"

$ws.Range("A7").Value = "OOOOOO"
$ws.Range("B7").Value = "6 Upper case O"

$ws.Range("A8").Value = "1111111111111111111111111111111111111111111111111111111111111111111111111111111111111111111111"
$ws.Range("B8").Value = "Number 1s"

$ws.Range("A9").Value = "£!_+@~#?|*"
$ws.Range("B9").Value = "Special Chars"

# --- Update selection to reflect where the author ended up ---
$ws.Range("A8").Select()
